# Applies the "Update with TFS items" changes to IPC.docx:
#   1. Fill in the first empty "TFS wit" value cell with "24451".
#   2. Fill in the second empty "TFS wit" value cell with "24451" and move
#      the trailing "_GoBack" bookmark there (right after the new text).
#   3. Merge the split "Project changeset version: " / "37418" runs into one run.
#   4. Merge the split "Starting point: " / "buffer management" runs into one run.
#   5. Merge the split "appl_" / "resb" / ".c" runs into a single "appl_resb.c" run.
#   6. As a consequence of moving the bookmark in step 2, the trailing paragraph
#      that used to hold "_GoBack" becomes a plain empty paragraph.

$d = $word.ActiveDocument

# --- 1. First "TFS wit" table: empty value cell -> "24451" ------------------
$tbl1 = $d.Tables(1)
$cell1 = $tbl1.Rows(4).Cells(2)
$cell1.Range.Text = "24451"

# --- 2. Second "TFS wit" table: empty value cell -> "24451" + _GoBack -------
$tbl1b = $d.Tables(1)
$cell2 = $tbl1b.Rows(9).Cells(2)

# Type "24451" plus a temporary placeholder character. Adding the placeholder
# means the insertion point for the bookmark (right after "24451") sits in the
# middle of a run instead of exactly on the paragraph-mark boundary, which is
# the position the COM layer handles correctly.
$cell2.Range.Text = "24451X"

$tbl1c = $d.Tables(1)
$cell2b = $tbl1c.Rows(9).Cells(2)
$afterDigits = $cell2b.Range.Start + 5

# Remove the old "_GoBack" bookmark (currently on the trailing paragraph at
# the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $bmOld = $d.Bookmarks("_GoBack")
    $bmOld.Delete()
}

# Re-create "_GoBack" collapsed right after the "24451" text we just typed.
$bmRange = $d.Range($afterDigits, $afterDigits)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary placeholder character now that the bookmark is anchored.
$placeholder = $d.Range($afterDigits, $afterDigits + 1)
$placeholder.Text = ""

# --- 3. Merge "Project changeset version: " + "37418" into one run ----------
$d.Content.Find.Execute("Project changeset version: 37418", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Project changeset version: 37418", 2) | Out-Null

# --- 4. Merge "Starting point: " + "buffer management" into one run ---------
$d.Content.Find.Execute("Starting point: buffer management", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Starting point: buffer management", 2) | Out-Null

# --- 5. Merge "appl_" + "resb" + ".c" into a single "appl_resb.c" run -------
$d.Content.Find.Execute("appl_resb.c", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "appl_resb.c", 2) | Out-Null
